# Gates Demo Final
#
# Applies the changes described by the commit to scan_childvacc_825_pg3.xlsx:
#  - survey!F60,F64,F68,F72,F76,F80,F84,F88 get a TRUE (hideInContents) flag
#  - survey's column E width grows (36 -> ~50.83 chars)
#  - choices!B2:C7 change from numeric 1 to the text "yes"
#  - the "choices" sheet becomes the active/selected sheet & tab,
#    with C8 selected there (the survey sheet view no longer owns the selection)

$wb = $excel.ActiveWorkbook

# --- survey sheet: add the new "hideInContents" TRUE markers -------------
$survey = $wb.Worksheets.Item("survey")

$survey.Range("F60").Value = $true
$survey.Range("F64").Value = $true
$survey.Range("F68").Value = $true
$survey.Range("F72").Value = $true
$survey.Range("F76").Value = $true
$survey.Range("F80").Value = $true
$survey.Range("F84").Value = $true
$survey.Range("F88").Value = $true

# Widen column E (closest the COM width<->pixel rounding lets us reach 50.83203125)
$survey.Columns.Item(5).ColumnWidth = 50

# --- choices sheet: data_value / display.text become "yes" ---------------
$choices = $wb.Worksheets.Item("choices")

$choices.Range("B2").Value = "yes"
$choices.Range("C2").Value = "yes"
$choices.Range("B3").Value = "yes"
$choices.Range("C3").Value = "yes"
$choices.Range("B4").Value = "yes"
$choices.Range("C4").Value = "yes"
$choices.Range("B5").Value = "yes"
$choices.Range("C5").Value = "yes"
$choices.Range("B6").Value = "yes"
$choices.Range("C6").Value = "yes"
$choices.Range("B7").Value = "yes"
$choices.Range("C7").Value = "yes"

# --- selection / active-tab bookkeeping -----------------------------------
# Move the live selection off survey (so it no longer carries tabSelected)
# and leave the choices sheet active with C8 selected, matching the diff.
$choices.Activate()
$choices.Range("C8").Select()
